$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.863.83'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.893.13'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.01'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.42'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '57.10'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +8.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.359'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0987'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.92'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +14.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.792'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.169.66'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.47%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.89%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.899.31'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.830.45'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.32'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0832'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '247.11'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.87%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.19'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.71'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +6.56%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.44'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.71'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.68%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.92%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.50%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.61%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -14.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.862'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0787'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +14.98%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.93'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.41%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Gas'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.29'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +27.43%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.99'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.316.91'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.35'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0810'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.67%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.75'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.94'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.56%  '
